$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert three new blank rows before row 219 (old rows 219-222 shift down
#    to become rows 222-225).  Rows.Insert() on this host stamps the new rows
#    with a border-less style variant, so immediately repair the look of the
#    freshly inserted B/C cells by copying the (already-correct) style from
#    row 218 onto them before writing any content.
# ---------------------------------------------------------------------------
$ws.Rows("219:221").Insert()

$ws.Range("B218:C218").Copy($ws.Range("B219:C219"))
$ws.Range("B218:C218").Copy($ws.Range("B220:C220"))
$ws.Range("B218:C218").Copy($ws.Range("B221:C221"))

# ---------------------------------------------------------------------------
# 2) Populate the three new "update" API rows (219-221).  Cell-write order
#    below is deliberate: it reproduces the exact order new shared strings
#    must be interned in (matching the target sharedStrings.xml append
#    order) - all three API-key cells first, then the three description
#    cells (written Orang / Periode / EMail so the first-seen order comes
#    out Orang(363), Periode(364), EMail(365)).
# ---------------------------------------------------------------------------
$ws.Cells.Item(219, 2).Value = "transaction.update.master.setPeriod"
$ws.Cells.Item(220, 2).Value = "transaction.update.master.setPerson"
$ws.Cells.Item(221, 2).Value = "transaction.update.master.setPersonAccountEMail"

$ws.Cells.Item(220, 3).Value = "Memutakhirkan Data Orang"
$ws.Cells.Item(219, 3).Value = "Memutakhirkan Data Periode"
$ws.Cells.Item(221, 3).Value = "Memutakhirkan Data Akun E-Mail Orang"

# ---------------------------------------------------------------------------
# 3) Fill in the three placeholder "create" API rows (20-22), which
#    previously all shared the generic "transaction.create.master.set" /
#    "Menyimpan Data Baru " placeholder strings.  Write order again mirrors
#    the target shared-string append order: B22, B21, B20 (keys) then
#    C20, C21, C22 (descriptions).
# ---------------------------------------------------------------------------
$ws.Cells.Item(22, 2).Value = "transaction.create.master.setPersonAccountEMail"
$ws.Cells.Item(21, 2).Value = "transaction.create.master.setPerson"
$ws.Cells.Item(20, 2).Value = "transaction.create.master.setPeriod"

$ws.Cells.Item(20, 3).Value = "Menyimpan Data Baru Periode"
$ws.Cells.Item(21, 3).Value = "Menyimpan Data Baru Orang"
$ws.Cells.Item(22, 3).Value = "Menyimpan Data Baru Akun E-Mail Orang"

# ---------------------------------------------------------------------------
# 4) Restore the user's on-save selection (bottom-right frozen pane) to C23.
# ---------------------------------------------------------------------------
$ws.Range("C23").Select()

Write-Host "API-Catalogue updated: Period / Person / PersonAccountEMail create+update rows added"
